# Adds the "Dependency injection" demo-project notes to the end of the
# document, reusing the existing trailing empty paragraph for the first
# new line and appending the rest as new paragraphs (mix of plain
# "Normal" headers and numbered-list ("List Paragraph", numId 1) bullets
# at two indent levels), finishing with a new trailing empty bullet
# paragraph at the deepest indent level.

$d = $word.ActiveDocument

function Add-PlainLine([string]$text) {
    $para = $d.Paragraphs.Last
    $r = $para.Range
    $r.InsertAfter($text)
    $r.LanguageID = "en-US"
    return $para
}

function Add-ListLine([string]$text, [int]$level) {
    $para = Add-PlainLine $text
    $para.Style = "List Paragraph"
    $para.Range.ListFormat.ListLevelNumber = $level
    return $para
}

function New-ParagraphAfterLast() {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
}

# 1) The document already ends with a single empty paragraph; turn it
#    into the "Dependency injection" heading (plain/Normal style).
Add-PlainLine "Dependency injection" | Out-Null
New-ParagraphAfterLast

# 2) Constructor injection (top-level bullet)
Add-ListLine "Constructor injection" 1 | Out-Null
New-ParagraphAfterLast

# 3) Use when you have require dependencies (sub bullet)
Add-ListLine "Use when you have require dependencies" 2 | Out-Null
New-ParagraphAfterLast

# 4) Generally recommended as first choice (sub bullet)
Add-ListLine "Generally recommended as first choice" 2 | Out-Null
New-ParagraphAfterLast

# 5) Setter injection (top-level bullet)
Add-ListLine "Setter injection" 1 | Out-Null
New-ParagraphAfterLast

# 6) Use this when you have optional dependencies (sub bullet)
Add-ListLine "Use this when you have optional dependencies" 2 | Out-Null
New-ParagraphAfterLast

# 7) If not provided, can provide reasonable default logic (sub bullet)
Add-ListLine "If not provided, can provide reasonable default logic" 2 | Out-Null
New-ParagraphAfterLast

# 8) Autowiring (top-level bullet)
Add-ListLine "Autowiring" 1 | Out-Null
New-ParagraphAfterLast

# 9) Spring will look for a class that matches (sub bullet)
Add-ListLine "Spring will look for a class that matches" 2 | Out-Null
New-ParagraphAfterLast

# 10) Matches by type and uses it automatically (sub bullet)
Add-ListLine "Matches by type and uses it automatically" 2 | Out-Null
New-ParagraphAfterLast

# 11) Final trailing paragraph: an *empty* sub-bullet at the same level.
#     Type placeholder text so the run/paragraph formatting (language +
#     list level) can be applied, then delete just the text, leaving the
#     paragraph-mark formatting (w:pPr/w:rPr) behind on the empty
#     paragraph - matching how Word represents an empty list item.
$finalPara = $d.Paragraphs.Last
$finalPara.Range.InsertAfter("TEMP")
$finalPara.Range.LanguageID = "en-US"
$finalPara.Style = "List Paragraph"
$finalPara.Range.ListFormat.ListLevelNumber = 2
$tempRange = $d.Range($finalPara.Range.Start, $finalPara.Range.Start + 4)
$tempRange.Delete()

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
